$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.047.96"
$ws.Range("E2").Value = "  -6.38%  "
$ws.Range("D3").Value = "2.188.99"
$ws.Range("E3").Value = "  -7.12%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.78"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.04"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.61%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -11.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.69"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0935"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -8.09%  "
$ws.Range("E13").Value = "  -4.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -9.34%  "
$ws.Range("D15").Value = "2.511.98"
$ws.Range("E15").Value = "  -7.30%  "
$ws.Range("E16").Value = "  -10.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -8.34%  "
$ws.Range("D18").Value = "2.187.84"
$ws.Range("E18").Value = "  -7.26%  "
$ws.Range("D19").Value = "41.002.65"
$ws.Range("E19").Value = "  -6.50%  "
$ws.Range("D20").Value = "0.0₃0936"
$ws.Range("E20").Value = "  -8.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.10"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.83%  "
$ws.Range("E22").Value = "  -7.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -8.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.07%  "
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("E28").Value = "  -5.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.12"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.46%  "
$ws.Range("E33").Value = "  -7.98%  "
$ws.Range("E34").Value = "  -6.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  -9.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.82"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +16.90%  "
$ws.Range("E39").Value = "  -6.76%  "
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -11.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.50"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.83"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -12.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  -5.16%  "
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0974"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -8.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.50"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.87%  "
$ws.Range("E51").Value = "  -6.23%  "
